# Update crypto price (D) and volume change (E) columns with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'64.297.81"
$ws.Cells.Item(2, 5).Value = "  +0.21%  "
$ws.Cells.Item(3, 4).Value = "'3.505.17"
$ws.Cells.Item(3, 5).Value = "  -0.50%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'589.67"
$ws.Cells.Item(5, 5).Value = "  +0.65%  "
$ws.Cells.Item(6, 4).Value = "'134.43"
$ws.Cells.Item(6, 5).Value = "  +0.20%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 4).Value = "'0.487"
$ws.Cells.Item(8, 5).Value = "  -0.52%  "
$ws.Cells.Item(9, 5).Value = "  +0.15%  "
$ws.Cells.Item(10, 5).Value = "  +2.54%  "
$ws.Cells.Item(11, 5).Value = "  +2.31%  "
$ws.Cells.Item(12, 4).Value = "'4.100.94"
$ws.Cells.Item(12, 5).Value = "  -0.57%  "
$ws.Cells.Item(13, 5).Value = "  +1.18%  "
$ws.Cells.Item(14, 4).Value = "'0.0000181"
$ws.Cells.Item(14, 5).Value = "  +1.13%  "
$ws.Cells.Item(15, 4).Value = "'3.502.62"
$ws.Cells.Item(15, 5).Value = "  -0.82%  "
$ws.Cells.Item(16, 4).Value = "'64.322.09"
$ws.Cells.Item(16, 5).Value = "  +0.18%  "
$ws.Cells.Item(17, 4).Value = "'25.68"
$ws.Cells.Item(17, 5).Value = "  -6.45%  "
$ws.Cells.Item(18, 5).Value = "  +0.66%  "
$ws.Cells.Item(19, 4).Value = "'5.75"
$ws.Cells.Item(19, 5).Value = "  +2.50%  "
$ws.Cells.Item(20, 4).Value = "'13.51"
$ws.Cells.Item(20, 5).Value = "  -2.86%  "
$ws.Cells.Item(21, 4).Value = "'393.08"
$ws.Cells.Item(21, 5).Value = "  +2.77%  "
$ws.Cells.Item(22, 5).Value = "  -0.16%  "
$ws.Cells.Item(23, 4).Value = "'3.644.59"
$ws.Cells.Item(23, 5).Value = "  -0.58%  "
$ws.Cells.Item(24, 4).Value = "'74.62"
$ws.Cells.Item(24, 5).Value = "  +0.77%  "
$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 5).Value = "  +0.02%  "
$ws.Cells.Item(26, 5).Value = "  +0.10%  "
$ws.Cells.Item(27, 4).Value = "'0.999"
$ws.Cells.Item(27, 5).Value = "  -0.05%  "
$ws.Cells.Item(28, 5).Value = "  -1.41%  "
$ws.Cells.Item(29, 5).Value = "  +1.03%  "
$ws.Cells.Item(30, 5).Value = "  -2.24%  "
$ws.Cells.Item(31, 5).Value = "  -7.62%  "
$ws.Cells.Item(32, 4).Value = "'3.527.34"
$ws.Cells.Item(32, 5).Value = "  -0.27%  "
$ws.Cells.Item(33, 4).Value = "'0.154"
$ws.Cells.Item(33, 5).Value = "  +5.44%  "
$ws.Cells.Item(34, 5).Value = "  +0.05%  "
$ws.Cells.Item(35, 4).Value = "'23.47"
$ws.Cells.Item(35, 5).Value = "  -0.44%  "
$ws.Cells.Item(36, 5).Value = "  -5.08%  "
$ws.Cells.Item(37, 4).Value = "'6.89"
$ws.Cells.Item(37, 5).Value = "  -0.94%  "
$ws.Cells.Item(38, 5).Value = "  -0.68%  "
$ws.Cells.Item(39, 4).Value = "'167.55"
$ws.Cells.Item(39, 5).Value = "  +4.22%  "
$ws.Cells.Item(40, 4).Value = "'0.0780"
$ws.Cells.Item(40, 5).Value = "  -0.78%  "
$ws.Cells.Item(41, 4).Value = "'0.812"
$ws.Cells.Item(41, 5).Value = "  -0.20%  "
$ws.Cells.Item(42, 5).Value = "  -0.04%  "
$ws.Cells.Item(43, 4).Value = "'25.10"
$ws.Cells.Item(43, 5).Value = "  -5.72%  "
$ws.Cells.Item(44, 4).Value = "'4.40"
$ws.Cells.Item(44, 5).Value = "  -0.30%  "
$ws.Cells.Item(45, 5).Value = "  +2.89%  "
$ws.Cells.Item(46, 5).Value = "  -3.97%  "
$ws.Cells.Item(47, 5).Value = "  -0.63%  "
$ws.Cells.Item(48, 4).Value = "'2.349.49"
$ws.Cells.Item(48, 5).Value = "  -5.37%  "
$ws.Cells.Item(49, 5).Value = "  -1.96%  "
$ws.Cells.Item(50, 5).Value = "  -1.55%  "
$ws.Cells.Item(51, 4).Value = "'21.09"
$ws.Cells.Item(51, 5).Value = "  -1.52%  "
